$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple text/value assignments (strings that Excel will not reinterpret as numbers)
$ws.Range("D2").Value = "60.760.14"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").Value = "2.997.16"
$ws.Range("E3").Value = "  -4.85%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -4.28%  "
$ws.Range("E6").Value = "  -5.73%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "2.994.37"
$ws.Range("E8").Value = "  -4.89%  "
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("E10").Value = "  -5.16%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -5.40%  "
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "3.484.15"
$ws.Range("E16").Value = "  -4.81%  "
$ws.Range("D17").Value = "60.988.40"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "2.981.90"
$ws.Range("E18").Value = "  -5.10%  "
$ws.Range("E19").Value = "  -4.85%  "
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("E21").Value = "  -5.00%  "
$ws.Range("E22").Value = "  -5.79%  "
$ws.Range("E23").Value = "  -5.63%  "
$ws.Range("E24").Value = "  -5.10%  "
$ws.Range("E25").Value = "  -5.74%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  -6.68%  "
$ws.Range("E29").Value = "  -7.08%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E30").Value = "  -6.00%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E31").Value = "  -9.03%  "
$ws.Range("E32").Value = "  -6.34%  "
$ws.Range("E33").Value = "  -8.90%  "
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("E35").Value = "  -8.14%  "
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").Value = "0.0₃0663"
$ws.Range("E38").Value = "  -5.41%  "
$ws.Range("E39").Value = "  -6.42%  "
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("E42").Value = "  -6.84%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.653.51"
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E44").Value = "  -7.12%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -6.00%  "
$ws.Range("E47").Value = "  -5.85%  "
$ws.Range("E48").Value = "  -7.04%  "
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("E50").Value = "  -3.73%  "
$ws.Range("E51").Value = "  -7.68%  "

# Cells whose new text looks like a number; force text storage, then restore original style
$c = $ws.Range("D4")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $origStyle

$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "562.98"
$c.Style = $origStyle

$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "128.08"
$c.Style = $origStyle

$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.24"
$c.Style = $origStyle

$c = $ws.Range("D12")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.430"
$c.Style = $origStyle

$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "32.73"
$c.Style = $origStyle

$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.119"
$c.Style = $origStyle

$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.23"
$c.Style = $origStyle

$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "437.92"
$c.Style = $origStyle

$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "13.13"
$c.Style = $origStyle

$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.663"
$c.Style = $origStyle

$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "78.85"
$c.Style = $origStyle

$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "12.53"
$c.Style = $origStyle

$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = $origStyle

$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.13"
$c.Style = $origStyle

$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.89"
$c.Style = $origStyle

$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.15"
$c.Style = $origStyle

$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "25.47"
$c.Style = $origStyle

$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0936"
$c.Style = $origStyle

$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.942"
$c.Style = $origStyle

$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "49.80"
$c.Style = $origStyle

$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0360"
$c.Style = $origStyle

$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.74"
$c.Style = $origStyle

$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.107"
$c.Style = $origStyle

$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "369.83"
$c.Style = $origStyle

$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = $origStyle

$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $origStyle

$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.236"
$c.Style = $origStyle

$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "118.08"
$c.Style = $origStyle

$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "32.94"
$c.Style = $origStyle

$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "23.35"
$c.Style = $origStyle
